$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GREEN = 5287936   # BGR for RGB 00B050 (accept)
$RED = 255         # BGR for RGB FF0000 (reject)

# Row 3: extra evaluations in F3 (Teamcity) and G3:H3 (SonarQube) - both "accepted"
$ws.Range("F3").Value = "a"
$ws.Range("F3").Font.Color = $GREEN

$ws.Range("G3").Value = "a"
$ws.Range("G3").MergeArea.Font.Color = $GREEN

# Remaining rows - fill in evaluation column C (merged C:D) with accept/reject
$answers = @{
    "C4"  = "a"
    "C5"  = "r"
    "C6"  = "r"
    "C7"  = "r"
    "C8"  = "a"
    "C9"  = "r"
    "C10" = "a"
    "C11" = "r"
    "C12" = "a"
    "C13" = "a"
    "C14" = "r"
    "C15" = "r"
    "C16" = "r"
    "C17" = "r"
}

foreach ($ref in $answers.Keys) {
    $val = $answers[$ref]
    $cell = $ws.Range($ref)
    $mergedRange = $cell.MergeArea
    $cell.Value = $val
    if ($val -eq "a") {
        $mergedRange.Font.Color = $GREEN
    } else {
        $mergedRange.Font.Color = $RED
    }
}

# Row 18 is special: base font for C18/D18 was plain Calibri (not the Marlett
# placeholder the other rows already had), so the color spreads across the
# merged C18:D18 range first, and only C18 is then switched to bold Marlett
# centered, matching the other answer cells.
$ws.Range("C18").Value = "r"
$ws.Range("C18:D18").Font.Color = $RED
$ws.Range("C18").Font.Name = "Marlett"
$ws.Range("C18").Font.Bold = $true
$ws.Range("C18").HorizontalAlignment = -4108

# Move the active selection like in the authored workbook
$ws.Range("F6").Select()
